$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Jengibre at "Vega Modelo de Temuco".
# It belongs right after the current row 66 (i.e. at row 67), pushing every
# row that used to be 67..139 down to 68..140.
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row 67 with the new record's data.
$ws.Cells.Item(67, 1).Value = 10
$ws.Cells.Item(67, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(67, 3).Value = "La Araucanía"
$ws.Cells.Item(67, 4).Value = 44629
$ws.Cells.Item(67, 5).Value = 9
$ws.Cells.Item(67, 6).Value = 100114007
$ws.Cells.Item(67, 7).Value = "Jengibre"
$ws.Cells.Item(67, 8).Value = "Sin especificar"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 15
$ws.Cells.Item(67, 11).Value = 25000
$ws.Cells.Item(67, 12).Value = 25000
$ws.Cells.Item(67, 13).Value = 25000
$ws.Cells.Item(67, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(67, 15).Value = "Perú"
$ws.Cells.Item(67, 16).Value = 1923
$ws.Cells.Item(67, 17).Value = 13
$ws.Cells.Item(67, 18).Value = "Hortaliza"
